$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = $ws.Range("B3").Value2
$ws.Range("D3").Value = $ws.Range("B3").Value2
$ws.Range("C9").Value = $ws.Range("B4").Value2
